$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First bullet under "Work on the Velocity Imaging Informatics System...":
#    The sentence itself is unchanged text-wise, but in the target OOXML the
#    five runs that spell it out get consolidated into a single run while the
#    following "trailing space" run (different rsidRPr) must stay untouched
#    and separate. A Find/Replace touching the paragraph normally sweeps that
#    trailing run in too (because it shares identical rPr), so we temporarily
#    bookmark it to break the adjacency before doing the replace, then drop
#    the bookmark again.
# ---------------------------------------------------------------------------
$sentence1 = "Working on developing and deploying cloud-based medical imaging services."

$r1 = $d.Content
$found1 = $r1.Find.Execute($sentence1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $guardStart = $r1.End
    $guardRange = $d.Range($guardStart, $guardStart + 1)
    $d.Bookmarks.Add("ZZ_GUARD1", $guardRange)
}

$r1b = $d.Content
$r1b.Find.Execute($sentence1, $true, $false, $false, $false, $false, $true, 1, $false, $sentence1, 2) | Out-Null

if ($d.Bookmarks.Exists("ZZ_GUARD1")) {
    $d.Bookmarks("ZZ_GUARD1").Delete()
}

# ---------------------------------------------------------------------------
# 2) Second bullet: "Working on visual and backend components of the
#    Velocity Imaging informatics software." becomes the new pipeline
#    sentence. This paragraph holds exactly one run, so a plain Find/Replace
#    is safe.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Working on visual and backend components of the Velocity Imaging informatics software.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Working on developing deployment pipelines for the cloud-based services.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Third bullet: "Working on algorithms to calculate-estimate radiation
#    doses." becomes the old second-bullet sentence. This paragraph's runs
#    span the whole paragraph (no trailing neighbour run), so a Find/Replace
#    that collapses them into one run matches the target exactly.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Working on algorithms to calculate-estimate radiation doses.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Working on visual and backend components of the Velocity Imaging informatics software.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Skills bullet (first one, ends in ", Python"): append new skills as a
#    brand-new run at the end of the paragraph, leaving the existing runs
#    untouched.
# ---------------------------------------------------------------------------
$skillsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.TrimEnd() -match ", Python$") {
        $skillsPara = $p
        break
    }
}
if ($skillsPara -ne $null) {
    $appendText1 = ", Terraform, AWS/ Azure Services"
    $paraEnd = $skillsPara.Range.End
    $skillsPara.Range.InsertAfter($appendText1)
    $newStart = $paraEnd - 1
    $newEnd = $newStart + $appendText1.Length
    $newRun = $d.Range($newStart, $newEnd)
    $newRun.Font.Name = "Helvetica"
}

# ---------------------------------------------------------------------------
# 5) Software/Frameworks bullet: ": TensorFlow" becomes ": TensorFlow,
#    OpenCV, QT, AWS/Azure Services, Docker, Kubernetes, Terraform" - i.e.
#    the existing "OpenCV, QT, AWS/Azure Services, Docker Container,
#    Kubernetes" tail (currently positioned after the _GoBack bookmark) is
#    logically moved ahead of the bookmark, "Docker Container" becomes
#    "Docker", and ", Terraform" is appended at the very end (after the
#    bookmark, matching the diff).
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute(", OpenCV, QT, AWS/Azure Services, Docker Container, Kubernetes", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$r5del = $d.Content
$found5 = $r5del.Find.Execute(", OpenCV, QT, AWS/Azure Services, Docker Container, Kubernetes", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$r6 = $d.Content
$r6.Find.Execute(": TensorFlow", $true, $false, $false, $false, $false, $true, 1, $false, `
    ": TensorFlow, OpenCV, QT, AWS/Azure Services, Docker, Kubernetes", 2) | Out-Null

$r7 = $d.Content
$found7 = $r7.Find.Execute("_GoBack_MARKER_NEVER_MATCHES_", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
